# Weekly update: insert a new data row for the latest reporting date right
# after the last "locked" historical row (row 100), pushing the existing
# rows 101-389 down by one (to 102-390).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 101; everything below (101-389) shifts to 102-390.
$ws.Rows.Item(101).Insert()

# The new row needs the same fixed/static columns (Mercado, Region, Codreg,
# Categoria, Variedad, Calidad, Unidad, Origen, Kg o Unidades, Clasificacion)
# and number formatting as every other row in this block. Copy row 102 (the
# row that just shifted down, still carrying the original formatting/values)
# into the new row 101, then overwrite the week-specific figures.
$ws.Range("A102:R102").Copy()
$ws.Range("A101").PasteSpecial(-4104) # xlPasteAll
$excel.CutCopyMode = $false

# New week's figures for row 101.
$ws.Cells.Item(101, 4).Value = 44914   # Fecha
$ws.Cells.Item(101, 10).Value = 129    # Volumen
$ws.Cells.Item(101, 11).Value = 1500   # Precio minimo
$ws.Cells.Item(101, 12).Value = 1800   # Precio maximo
$ws.Cells.Item(101, 13).Value = 1649   # Precio promedio ponderado
$ws.Cells.Item(101, 16).Value = 550    # Precio $/Kg
